$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.555.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "'3.484.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'581.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'147.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("D7").Value = "'3.484.37"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("D10").Value = "'7.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'0.126"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "'0.407"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.54%  "
$ws.Range("D13").Value = "'4.073.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "'29.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.44%  "
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "'3.481.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "'63.397.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "'6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.71%  "
$ws.Range("D20").Value = "'14.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.35%  "
$ws.Range("D21").Value = "'9.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").Value = "'391.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("D24").Value = "'75.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'3.608.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "'0.181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").Value = "'7.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").Value = "'2.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("D35").Value = "'23.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").Value = "'7.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'32.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.98%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").Value = "'1.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.70%  "
$ws.Range("D40").Value = "'170.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "'3.515.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "'0.0767"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").Value = "'0.801"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("D47").Value = "'4.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'2.628.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("D49").Value = "'2.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.53%  "
$ws.Range("D50").Value = "'23.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("E51").Value = "  +2.43%  "
